$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Cade Cunningham"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Detroit Pistons"

$ws.Range("A4").Value = "Cam Thomas"
$ws.Range("B4").Value = "SG,SF"
$ws.Range("C4").Value = "Brooklyn Nets"

$ws.Range("A5").Value = "Brandon Miller"
$ws.Range("B5").Value = "SG,SF"
$ws.Range("C5").Value = "Charlotte Hornets"

$ws.Range("A6").Value = "Royce O'Neale"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Phoenix Suns"

$ws.Range("A8").Value = "Julius Randle"
$ws.Range("B8").Value = "PF"
$ws.Range("C8").Value = "Minnesota Timberwolves"

$ws.Range("A9").Value = "Cameron Johnson"
$ws.Range("B9").Value = "SF,PF"
$ws.Range("C9").Value = "Brooklyn Nets"

$ws.Range("A11").Value = "Bam Adebayo"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Miami Heat"

$ws.Range("A12").Value = "Brandon Ingram"
$ws.Range("B12").Value = "SG,SF,PF"
$ws.Range("C12").Value = "New Orleans Pelicans"

$ws.Range("A13").Value = "Yves Missi"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "New Orleans Pelicans"

$ws.Range("A14").Value = "Derrick White"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Boston Celtics"

$ws.Range("A15").Value = "Gradey Dick"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Toronto Raptors"

$ws.Range("A16").Value = "Anthony Davis"
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Los Angeles Lakers"
